# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns,
# plus the PancakeSwap/LEO row swap (rows 26/27).
#
# Helper: some new Price values parse as "clean" decimals (e.g. "300.00",
# "17.78") which this host would otherwise auto-coerce to numeric cells.
# The source workbook stores every Price/Volume cell as text, so for those
# values we briefly mark the cell as text (NumberFormat "@"), assign the
# literal string, then restore the cell's style to Normal so no stray
# number-format sticks around on the cell.
function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26 / 27 swap (PancakeSwap <-> LEO) ---
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D26") "4.01"
$ws.Range("E26").Value = "  -0.36%  "

$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D27") "2.41"
$ws.Range("E27").Value = "  -1.49%  "

# --- Price (D) and Volume(1h) (E) updates for other rows ---

$ws.Range("D2").Value = "42.770.51"
$ws.Range("E2").Value = "  -0.71%  "

$ws.Range("D3").Value = "2.288.63"
$ws.Range("E3").Value = "  -1.05%  "

$ws.Range("E4").Value = "  +0.01%  "

Set-TextValue $ws.Range("D5") "300.00"
$ws.Range("E5").Value = "  -0.77%  "

Set-TextValue $ws.Range("D6") "96.10"
$ws.Range("E6").Value = "  -3.01%  "

$ws.Range("E7").Value = "  +0.55%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  -3.99%  "

$ws.Range("E10").Value = "  -0.95%  "

$ws.Range("E11").Value = "  -0.53%  "

$ws.Range("E12").Value = "  +0.58%  "

Set-TextValue $ws.Range("D13") "17.78"
$ws.Range("E13").Value = "  -1.13%  "

Set-TextValue $ws.Range("D14") "6.73"
$ws.Range("E14").Value = "  -2.36%  "

$ws.Range("D15").Value = "2.645.52"
$ws.Range("E15").Value = "  -0.92%  "

$ws.Range("D16").Value = "2.293.18"
$ws.Range("E16").Value = "  -3.83%  "

$ws.Range("E17").Value = "  -2.46%  "

$ws.Range("D18").Value = "42.675.13"
$ws.Range("E18").Value = "  -0.69%  "

Set-TextValue $ws.Range("D19") "12.75"
$ws.Range("E19").Value = "  -5.41%  "

$subscript3 = [char]0x2083
Set-TextValue $ws.Range("D20") "0.0$($subscript3)0903"
$ws.Range("E20").Value = "  -1.08%  "

Set-TextValue $ws.Range("D21") "6.02"
$ws.Range("E21").Value = "  -2.66%  "

Set-TextValue $ws.Range("D22") "67.60"
$ws.Range("E22").Value = "  -0.52%  "

Set-TextValue $ws.Range("D23") "239.93"
$ws.Range("E23").Value = "  -0.23%  "

Set-TextValue $ws.Range("D24") "2.13"
$ws.Range("E24").Value = "  -1.67%  "

$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("E28").Value = "  -0.11%  "

Set-TextValue $ws.Range("D29") "165.68"
$ws.Range("E29").Value = "  -2.16%  "

Set-TextValue $ws.Range("D30") "2.02"
$ws.Range("E30").Value = "  -1.71%  "

$ws.Range("E31").Value = "  -2.03%  "

Set-TextValue $ws.Range("D32") "32.69"
$ws.Range("E32").Value = "  -2.58%  "

$ws.Range("E33").Value = "  +0.10%  "

$ws.Range("E34").Value = "  -2.04%  "

$ws.Range("E35").Value = "  -4.61%  "

Set-TextValue $ws.Range("D36") "16.90"
$ws.Range("E36").Value = "  -7.94%  "

$ws.Range("E37").Value = "  -1.72%  "

Set-TextValue $ws.Range("D38") "0.0683"
$ws.Range("E38").Value = "  -1.45%  "

$ws.Range("E39").Value = "  -1.56%  "

$ws.Range("E40").Value = "  -3.41%  "

$ws.Range("E41").Value = "  -0.33%  "

$ws.Range("E42").Value = "  -2.09%  "

$ws.Range("D43").Value = "2.013.69"
$ws.Range("E43").Value = "  +0.90%  "

$ws.Range("E44").Value = "  -3.07%  "

Set-TextValue $ws.Range("D45") "10.06"
$ws.Range("E45").Value = "  -0.32%  "

Set-TextValue $ws.Range("D46") "2.10"
$ws.Range("E46").Value = "  -3.23%  "

Set-TextValue $ws.Range("D47") "17.10"
$ws.Range("E47").Value = "  -2.62%  "

Set-TextValue $ws.Range("D48") "2.77"
$ws.Range("E48").Value = "  -2.67%  "

Set-TextValue $ws.Range("D49") "2.92"
$ws.Range("E49").Value = "  -2.64%  "

$ws.Range("D50").Value = "2.514.69"
$ws.Range("E50").Value = "  -0.85%  "

Set-TextValue $ws.Range("D51") "52.99"
$ws.Range("E51").Value = "  -3.21%  "
